$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 866-867, shifting existing rows 866:925 down to 868:927
$ws.Range("A866:R867").EntireRow.Insert()

# Row 866: new "Primera" quality record dated 2023-01-05 (serial 44931)
$ws.Range("A866").Value = 3
$ws.Range("B866").Value = "Femacal de La Calera"
$ws.Range("C866").Value = "Coquimbo"
$ws.Range("D866").Value = 44931
$ws.Range("D866").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E866").Value = 5
$ws.Range("F866").Value = 100112023
$ws.Range("G866").Value = "Brócoli"
$ws.Range("H866").Value = "Sin especificar"
$ws.Range("I866").Value = "Primera"
$ws.Range("J866").Value = 1500
$ws.Range("K866").Value = 800
$ws.Range("L866").Value = 800
$ws.Range("M866").Value = 800
$ws.Range("N866").Value = "$/unidad"
$ws.Range("O866").Value = "Provincia de Quillota"
$ws.Range("P866").Value = 800
$ws.Range("Q866").Value = 1
$ws.Range("R866").Value = "Hortaliza"

# Row 867: new "Segunda" quality record, same date
$ws.Range("A867").Value = 3
$ws.Range("B867").Value = "Femacal de La Calera"
$ws.Range("C867").Value = "Coquimbo"
$ws.Range("D867").Value = 44931
$ws.Range("D867").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E867").Value = 5
$ws.Range("F867").Value = 100112023
$ws.Range("G867").Value = "Brócoli"
$ws.Range("H867").Value = "Sin especificar"
$ws.Range("I867").Value = "Segunda"
$ws.Range("J867").Value = 1600
$ws.Range("K867").Value = 700
$ws.Range("L867").Value = 700
$ws.Range("M867").Value = 700
$ws.Range("N867").Value = "$/unidad"
$ws.Range("O867").Value = "Provincia de Quillota"
$ws.Range("P867").Value = 700
$ws.Range("Q867").Value = 1
$ws.Range("R867").Value = "Hortaliza"
